# Add three new "satellite / landsat" vegetation-index rows to the env sheet
# (ndvi, ndmi, msavi), matching commit "added indices to `env`".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 67 : ndvi -------------------------------------------------------
$ws.Range("D67").ClearFormats()

$ws.Range("A67").Value = "satellite"
$ws.Range("B67").Value = "landsat"
$ws.Range("C67").Value = "AusCover"
$ws.Range("D67").Value = "surface reflectance"
$ws.Range("M67").Value = "normalised difference vegetation index"
$ws.Range("E67").Value = "ndvi"
$ws.Range("N67").Value = "(nir - red) / (nir + red)"
$ws.Range("O67").Value = "photosynthetic activity"
$ws.Range("F67").Value = 1
$ws.Range("G67").Value = "mean"
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = "each"
$ws.Range("J67").Value = 1989
$ws.Range("J67").NumberFormat = "0"
$ws.Range("K67").Value = 10000
$ws.Range("L67").Value = "units"

# --- Row 68 : ndmi ---------------------------------------------------------
$ws.Range("A68").Value = "satellite"
$ws.Range("B68").Value = "landsat"
$ws.Range("C68").Value = "AusCover"
$ws.Range("D68").Value = "surface reflectance"
$ws.Range("E68").Value = "ndmi"
$ws.Range("M68").Value = "normalised difference moisture index"
$ws.Range("N68").Value = "(nir - swir1) / (nir + swir1)"
$ws.Range("O68").Value = "vegetation water content"
$ws.Range("F68").Value = 1
$ws.Range("G68").Value = "mean"
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = "each"
$ws.Range("J68").Value = 1989
$ws.Range("J68").NumberFormat = "0"
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = "units"

# --- Row 69 : msavi ---------------------------------------------------------
$ws.Range("A69").Value = "satellite"
$ws.Range("B69").Value = "landsat"
$ws.Range("C69").Value = "AusCover"
$ws.Range("D69").Value = "surface reflectance"
$ws.Range("E69").Value = "msavi"
$ws.Range("M69").Value = "modified soil adjusted vegetation index"
$ws.Range("N69").Value = "(2 * nir + 1 - sqrt((2 * nir + 1)^2 - 8 * (nir - red)))/2"
$ws.Range("O69").Value = "minimizes the effect of bare soil on the Soil Adjusted Vegetation Index (SAVI)"
$ws.Range("F69").Value = 1
$ws.Range("G69").Value = "mean"
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = "each"
$ws.Range("J69").Value = 1989
$ws.Range("J69").NumberFormat = "0"
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = "units"

# --- Clean up the now-unused "Hyperlink" cell style left over from the
#     stray D67 cell that used to hold a (now removed) hyperlink. ---------
$wb.Styles("Hyperlink").Delete() | Out-Null

# --- Restore the recorded selection state --------------------------------
$ws.Range("J61").Select() | Out-Null
